$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.850.42"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "2.588.46"
$ws.Range("E3").Value = "  +1.58%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.39"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.77"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.595"
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.28"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "3.057.95"
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("D15").Value = "62.840.19"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("D17").Value = "2.599.87"
$ws.Range("E17").Value = "  +2.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.29"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "342.19"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("E20").Value = "  +1.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.69"
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.29"
$ws.Range("E23").Value = "  +2.07%  "
$ws.Range("D24").Value = "2.715.65"
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.59"
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -2.67%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.85"
$ws.Range("E28").Value = "  +2.06%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.31"
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("E30").Value = "  -2.33%  "
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "471.67"
$ws.Range("E32").Value = "  +13.23%  "
$ws.Range("D33").Value = "0.0₃0822"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "176.92"
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("E35").Value = "  +4.26%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.406"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.00"
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.51"
$ws.Range("E39").Value = "  +3.86%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.70"
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "158.58"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.75"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.14"
$ws.Range("E44").Value = "  +1.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.632"
$ws.Range("E45").Value = "  +4.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0542"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0967"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0236"
$ws.Range("E48").Value = "  -0.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.41"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.72"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("E51").Value = "  +1.08%  "
